$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.059.30'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '3.134.24'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  -0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '570.79'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.28%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '161.63'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -3.94%  '
$ws.Range('E7').Value = '  -0.21%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.573'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -5.93%  '
$ws.Range('D9').Value = '3.147.10'
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('E10').Value = '  -3.31%  '
$ws.Range('E11').Value = '  -2.90%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.384'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').Value = '3.686.25'
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').Value = '64.138.83'
$ws.Range('E15').Value = '  -0.41%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '25.02'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '3.141.66'
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('E18').Value = '  -3.30%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '400.70'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -4.33%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '5.24'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '12.52'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -2.45%  '
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('E23').Value = '  +3.15%  '
$ws.Range('E24').Value = '  +0.05%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '67.96'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -2.53%  '
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('E27').Value = '  -5.37%  '
$ws.Range('E28').Value = '  -5.04%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '8.75'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.22%  '
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('E31').Value = '  +0.06%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.80'
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '21.10'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -2.60%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '159.38'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('E35').Value = '  -1.08%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '4.80'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -4.24%  '
$ws.Range('E37').Value = '  -2.27%  '
$ws.Range('E38').Value = '  -1.47%  '
$ws.Range('D39').Value = '2.664.18'
$ws.Range('E39').Value = '  -2.14%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.67'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -1.90%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '23.62'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -2.79%  '
$ws.Range('E42').Value = '  -2.45%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '38.29'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -2.16%  '
$ws.Range('E44').Value = '  -2.72%  '
$ws.Range('E45').Value = '  -1.10%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '5.42'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -2.96%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0254'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -2.58%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '286.37'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -2.34%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '21.02'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -2.43%  '
$ws.Range('E50').Value = '  -0.33%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0974'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.20%  '
